# Updates the cryptocurrency price/volume table (and fixes a row-order
# swap between ImmutableX and WrappedliquidstakedEther2.0) to match the
# latest scrape, per the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.909.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.758.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3763"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3359"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.89"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.126"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07168"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.36"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.177"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.166"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.762.07"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001050"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06569"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.29"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.87"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.258"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.934.98"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.66"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -9.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.395"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.51"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.69"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.319"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -10.38%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.963.11"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.269"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -17.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "130.92"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.015"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.781"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08755"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.14"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -9.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02337"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6553"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06178"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.109"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2102"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.207"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.459"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.023"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.64"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.26%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6004"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.997"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07208"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.174"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.18%  "
